# chore: add monthly employment outputs
#
# Updates the 피보험자수(insured persons) figures for 경기도 시군 (Gyeonggi
# municipalities) to reflect the newly published monthly employment output,
# and adds 광명시 as a newly flagged "주의"(caution) region in the
# 3-month-consecutive 시군 sheet (which pushes 연천군/파주시 down one row).

$wb = $excel.ActiveWorkbook

# --- Sheet "요약_권역별": 시군 / 피보험자수 row (row 7) ---
$ws1 = $wb.Worksheets.Item("요약_권역별")
$ws1.Range("C7").Value = 27
$ws1.Range("E7").Value = 3

# --- Sheet "요약_전월대비": gyeonggi_city / 피보험자수 row (row 7) ---
$ws2 = $wb.Worksheets.Item("요약_전월대비")
$ws2.Range("C7").Value = 27
$ws2.Range("E7").Value = 3

# --- Sheet "3개월연속_시군": insert a new 피보험자수/광명시 "주의" row ---
# Before: row15=연천군(주의), row16=파주시(주의)
# After : row15=광명시(주의, new), row16=연천군(주의, shifted), row17=파주시(주의, shifted)
$ws4 = $wb.Worksheets.Item("3개월연속_시군")
$ws4.Rows.Item(16).Insert()

$ws4.Range("A15").Value = "피보험자수"
$ws4.Range("B15").Value = "광명시"
$ws4.Range("C15").Value = "주의"
$ws4.Range("D15").Value = "주의"
$ws4.Range("E15").Value = "주의"

$ws4.Range("A16").Value = "피보험자수"
$ws4.Range("B16").Value = "연천군"
$ws4.Range("C16").Value = "주의"
$ws4.Range("D16").Value = "주의"
$ws4.Range("E16").Value = "주의"

# --- Sheet "주요지역_시군": 상실자수 / 양평군 row (row 32) ---
$ws6 = $wb.Worksheets.Item("주요지역_시군")
$ws6.Range("C32").Value = 3
$ws6.Range("E32").Value = 3
